# Update the two remaining LR-pair rows with newly recomputed TPM-derived
# values, then drop the old third data row (ECs -> ECs) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs (sending) -> Efnb3/Rhbdl2 -> ECs (target) -------------
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Rhbdl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1498043333333333
$ws.Range("H2").Value = 0.449413
$ws.Range("I2").Value = 0.08722868471333377
$ws.Range("J2").Value = 0.08722868471333377
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.415892
$ws.Range("N2").Value = 1.247676
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.06230242379866667
$ws.Range("R2").Value = 0.560721814188
$ws.Range("S2").Value = 0.08722868471333377
$ws.Range("T2").Value = 0.08722868471333377

# --- Row 3: MuSCs (sending) -> Efnb3/Rhbdl2 -> ECs (target) ------------
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Rhbdl2"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.567570333333333
$ws.Range("H3").Value = 4.702711
$ws.Range("I3").Value = 0.9127713152866662
$ws.Range("J3").Value = 0.9127713152866662
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.415892
$ws.Range("N3").Value = 1.247676
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.6519399610706667
$ws.Range("R3").Value = 5.867459649636
$ws.Range("S3").Value = 0.9127713152866662
$ws.Range("T3").Value = 0.9127713152866662

# --- Drop the old row 4 (ECs -> ECs), which no longer exists -----------
$ws.Rows("4:4").Delete()
